$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the end time recorded on row 80 (2014-03-20): 12:15 -> 12:30.
$ws.Range("E80").Value = 0.52083333333333337

# Move the active selection to E81, matching the saved workbook state.
$ws.Range("E81").Select()
